# Insert one new data row at row 152 (pushing the existing rows 152-195
# down to 153-196) and populate it with the new price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(152).Insert()

$ws.Range("A152").Value = 10
$ws.Range("B152").Value = "Vega Modelo de Temuco"
$ws.Range("C152").Value = "La Araucanía"
$ws.Range("D152").Value = 45135
$ws.Range("E152").Value = 9
$ws.Range("F152").Value = 100112035
$ws.Range("G152").Value = "Bruselas (repollito)"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 30
$ws.Range("K152").Value = 25000
$ws.Range("L152").Value = 25000
$ws.Range("M152").Value = 25000
$ws.Range("N152").Value = "$/malla 15 kilos"
$ws.Range("O152").Value = "Provincia de Quillota"
$ws.Range("P152").Value = 1667
$ws.Range("Q152").Value = 15
$ws.Range("R152").Value = "Hortaliza"
